# Updated legacy GSC export data.
#
# The "Chart" sheet holds a rolling window of daily Search Console metrics
# (Date / Not indexed / Indexed / Impressions). The export has rolled
# forward by one day: the oldest date row is dropped and every remaining
# row shifts up, which also drops the now out-of-window last row.
#
# Deleting row 2 (the first data row, right under the header) shifts all
# subsequent rows up by one and shrinks the used range accordingly - this
# reproduces the export refresh exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows.Item(2).Delete()
